$wb = $excel.ActiveWorkbook

# --- Belgium sheet is no longer the active/selected tab; its selection resets to the full used range ---
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Activate()
$belgium.Range("A1:D23").Select()

# --- Create the new "Czech" sheet as a copy of "Belgium" (keeps styles/merges/etc) ---
$belgium.Copy([System.Type]::Missing, $belgium)
$czech = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Name = "Czech"

# Czech market only has 2 repeaters rows (Belgium's rows 20 "ZXF" and 21 "ZXFEV" do not apply)
$czech.Rows("20:21").Delete()

# Update the market name / user story cells for the Czech sheet
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1735/T1741"

# Column widths differ on the Czech sheet (closest values reachable through this host)
$czech.Columns("B").ColumnWidth = 26.3
$czech.Columns("C").ColumnWidth = 16.15
$czech.Columns("D").ColumnWidth = 16.65

# Czech is now the selected/active sheet with the same B4 selection Belgium used to have
$czech.Activate()
$czech.Range("B4").Select()
